# Update "想去人数" (number of people interested) counts for two events
# that appear on both the "展览" sheet and the combined "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# "展览" sheet: rows 4 and 6 hold the two affected events.
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 836
$wsExpo.Range("F6").Value = 28

# "全部类型" sheet: same two events appear at rows 5 and 7.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 836
$wsAll.Range("F7").Value = 28
